$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells already holding plain-text numeric-looking values,
# used as a Copy() source so new cells keep text type (t="s") instead
# of Excel auto-coercing strings like "1"/"3" into numbers on assignment.
$srcOne = $ws.Cells.Item(18, 2)   # "1" as text
$srcThree = $ws.Cells.Item(45, 2) # "3" as text

function Set-TextCell($cell, $text) {
    if ($text -eq "1") {
        $srcOne.Copy($cell)
    } elseif ($text -eq "3") {
        $srcThree.Copy($cell)
    } else {
        $cell.Value = $text
    }
}

$data = @(
    ,@(109, "14-12-2021 11:16", "Hola", $null)
    ,@(110, "14-12-2021 11:16", "Menu", $null)
    ,@(111, "14-12-2021 11:16", "menu", "STEP_2_1")
    ,@(112, "14-12-2021 11:17", "3", $null)
    ,@(113, "14-12-2021 11:17", "A1", $null)
    ,@(114, "14-12-2021 11:17", "J1", $null)
    ,@(115, "14-12-2021 11:19", "Hola", $null)
    ,@(116, "14-12-2021 11:19", "Menu", $null)
    ,@(117, "14-12-2021 11:19", "menu", "STEP_2_1")
    ,@(118, "14-12-2021 11:19", "3", $null)
    ,@(119, "14-12-2021 11:19", "A1", $null)
    ,@(120, "14-12-2021 11:28", "Hola", $null)
    ,@(121, "14-12-2021 11:28", "A1", $null)
    ,@(122, "14-12-2021 11:49", "A1", $null)
    ,@(123, "14-12-2021 11:53", "A1", $null)
    ,@(124, "14-12-2021 11:53", "Hola", $null)
    ,@(125, "14-12-2021 11:53", "Menu", $null)
    ,@(126, "14-12-2021 11:53", "menu", "STEP_2_1")
    ,@(127, "14-12-2021 11:53", "3", $null)
    ,@(128, "14-12-2021 11:54", "Hola", $null)
    ,@(129, "14-12-2021 11:54", "Menu", $null)
    ,@(130, "14-12-2021 11:54", "menu", "STEP_2_1")
    ,@(131, "14-12-2021 11:54", "1", $null)
    ,@(132, "14-12-2021 11:54", "J3", $null)
)

foreach ($row in $data) {
    $r = $row[0]
    Set-TextCell $ws.Cells.Item($r, 1) $row[1]
    Set-TextCell $ws.Cells.Item($r, 2) $row[2]
    if ($row[3] -ne $null) {
        Set-TextCell $ws.Cells.Item($r, 3) $row[3]
    }
}
